# Daily attendance processing - 2025-10-06 19:15:21
# Reorders the comma-separated "Recorded By" names in column G for
# specific rows on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "System, system, backup@backdoor.com"
    3   = "dnasr281@gmail.com, System"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    6   = "dnasr281@gmail.com, System"
    11  = "dnasr281@gmail.com, System"
    12  = "dnasr281@gmail.com, System"
    13  = "dnasr281@gmail.com, System"
    29  = "System, system, backup@backdoor.com"
    30  = "dnasr281@gmail.com, System"
    32  = "System, backup@backdoor.com"
    33  = "dnasr281@gmail.com, System"
    38  = "dnasr281@gmail.com, System"
    39  = "dnasr281@gmail.com, System"
    40  = "dnasr281@gmail.com, System"
    56  = "System, system, backup@backdoor.com"
    57  = "dnasr281@gmail.com, System"
    58  = "System, backup@backdoor.com"
    59  = "System, backup@backdoor.com"
    60  = "dnasr281@gmail.com, System"
    65  = "dnasr281@gmail.com, System"
    66  = "dnasr281@gmail.com, System"
    67  = "dnasr281@gmail.com, System"
    84  = "System, backup@backdoor.com"
    85  = "System, backup@backdoor.com"
    89  = "dnasr281@gmail.com, System"
    90  = "dnasr281@gmail.com, admin@admin.com"
    93  = "dnasr281@gmail.com, System"
    110 = "System, backup@backdoor.com"
    111 = "System, backup@backdoor.com"
    115 = "dnasr281@gmail.com, System"
    116 = "dnasr281@gmail.com, admin@admin.com"
    119 = "dnasr281@gmail.com, System"
    136 = "System, backup@backdoor.com"
    137 = "System, backup@backdoor.com"
    141 = "dnasr281@gmail.com, System"
    142 = "dnasr281@gmail.com, admin@admin.com"
    145 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
